$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was recorded for Papaya (Vega Modelo de Temuco).
# It belongs chronologically before the existing row 13, so insert a fresh row
# there; Excel shifts every row from 13..85 down to 14..86 automatically.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new observation.
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "Vega Modelo de Temuco"
$ws.Range("C13").Value = "La Araucanía"
$ws.Range("D13").Value = 44971
$ws.Range("E13").Value = 9
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100108
$ws.Range("H13").Value = "Tropicales y subtropicales"
$ws.Range("I13").Value = 100108004
$ws.Range("J13").Value = "Papaya"
$ws.Range("K13").Value = "Cultivar IV Región"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 65
$ws.Range("N13").Value = 45000
$ws.Range("O13").Value = 45000
$ws.Range("P13").Value = 45000
$ws.Range("Q13").Value = "`$/caja 15 kilos granel"
$ws.Range("R13").Value = "Provincia del Elquí"
$ws.Range("S13").Value = 3000
$ws.Range("T13").Value = 15
